# "Added 2020 results + clean up"
#
# - Remove the (now unused) "Assumptions" sheet; this also drops the
#   orphaned "Super small # =" shared string and shifts the remaining
#   shared-string indices down by one (handled automatically by the
#   engine's shared-strings GC on save).
# - Protect the remaining data sheets (Election Results by State,
#   Uncontested Races, Uncontested by State PIVOT, EXPORT).
# - Leave EXPORT as the active/selected sheet (it is now the 4th tab
#   instead of the 5th after the Assumptions sheet is removed).

$wb = $excel.ActiveWorkbook

# Delete the "Assumptions" sheet entirely.
$assumptions = $wb.Worksheets.Item("Assumptions")
$assumptions.Delete()

# Protect the remaining sheets.
$wb.Worksheets.Item("Election Results by State").Protect()
$wb.Worksheets.Item("Uncontested Races").Protect()
$wb.Worksheets.Item("Uncontested by State PIVOT").Protect()
$wb.Worksheets.Item("EXPORT").Protect()

# Keep EXPORT as the active tab (now index 3 instead of 4).
$wb.Worksheets.Item("EXPORT").Activate()
